$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 (A39 = "Travis Winston") gets a value in column B, the same
# "BYAC290" course code already used for other rows in this column.
$ws.Range("B39").Value = "BYAC290"

# Reflect the author's scrolled/selected position in the sheet view:
# the window had scrolled further down and the active cell moved from
# B39 to B40.
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B40").Select()
